# ARKCORR-22 Added elvis operator to drools rules.
#
# The "CONDITION" column (C) of the "Set Due Date ... Queue" rows checks
# queue.name directly, which NPEs when a case file has no queue assigned.
# Switch each of those conditions to use the safe-navigation ("Elvis")
# operator: queue.name -> queue?.name

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C27").Value2 = 'queue?.name == "Intake"'
$ws.Range("C28").Value2 = 'queue?.name == "Fulfill"'
$ws.Range("C29").Value2 = 'queue?.name == "Supervisor Approval"'
$ws.Range("C30").Value2 = 'queue?.name == "Executive Approval"'
$ws.Range("C31").Value2 = 'queue?.name == "Release"'

# Reflect the editor's scroll position / current selection at save time.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 2
$ws.Range("C31").Select()
